$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: TestCase ID
$ws.Range("A2").Value = 'TC_14'

# B2: TestCaseDescription
$ws.Range("B2").Value = 'Verify that the ERROR message is received when missing HTTP Header X-CSR-SECURITY_TOKEN.'

# C2: ScenarioType -- unchanged by this commit ("Post" stays "Post"), left alone.

# D2: Input Request JSON
$ws.Range("D2").Value = '
{
	"meta":
	{
		"userName":"3791813"
	},
	"geopoliticalAffiliationType":
	{
		"affiliationTypeCode":"aaa",
		"affiliationTypeName":"Affilcodttt"
	}
}'

# E2: Input Data
$ws.Range("E2").Value = 'Input_UserName: 3791813
Input_affiliationTypeCode: aaa
Input_affiliationTypeName: Affilcodttt
'

# F2: DB Data -> now just "NA"
$ws.Range("F2").Value = 'NA'

# G2: WS Status
$ws.Range("G2").Value = 'HTTP/1.1 401 Unauthorized'

# H2: WS Status Code -- numeric-looking text "401"; must stay text (not auto-convert to a number).
# Build it on a scratch cell via a formula (forces a string result), then Copy/PasteSpecial values
# onto H2 so the literal shared-string is written without Excel coercing it to a Double and without
# disturbing H2's existing cell style (NumberFormat="@" on H2 directly would allocate a new style).
$scratch = $ws.Range("Z1")
$scratch.Formula = '="401"'
$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$scratch.Clear()

# I2: WS Response JSON
$ws.Range("I2").Value = '
{
	"meta":
	{
		"version":"1.0.0",
		"timeStamp":"2020-06-22T09:32:53.175+0000",
		"statusCode":"401",
		"message":
		{
			"status":"ERROR",
			"internalMessage":"Security Error",
			"data":
			{
				"errorMessage":"Exception occurred in Security Filter"
			}
		}
	},
	"errors":
	[
		{
			"fieldName":"NA",
			"message":"Missing HTTP header X-CSR-SECURITY_TOKEN"
		}
	]
}'

# J2 / K2 unchanged (TestResult "Pass", Comments blank) -- left as-is.

